$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (U2) - Mid X / Ref X / Pad X updates
$ws.Range("D3").Value = "7.747mm"
$ws.Range("F3").Value = "7.747mm"
$ws.Range("H3").Value = "9.017mm"

# Row 4 (PROG) - Mid/Ref/Pad X-Y and Rotation updates
$ws.Range("D4").Value = "-12.7mm"
$ws.Range("E4").Value = "-14.097mm"
$ws.Range("F4").Value = "-12.7mm"
$ws.Range("G4").Value = "-14.097mm"
$ws.Range("H4").Value = "-6.35mm"
$ws.Range("I4").Value = "-14.097mm"
$ws.Range("L4").Value = 180

# Row 5 - Designator renamed from CN1 to WIEGAND
$ws.Range("A5").Value = "WIEGAND"

# Row 7 (U3 / Level Shifter) - Footprint name + coordinates updates
$ws.Range("C7").Value = "Level Shifter Footprint"
$ws.Range("D7").Value = "-10.776mm"
$ws.Range("E7").Value = "0.951mm"
$ws.Range("F7").Value = "-17.145mm"
$ws.Range("G7").Value = "9.017mm"
$ws.Range("H7").Value = "-15.875mm"
$ws.Range("I7").Value = "7.301mm"
